$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 29 (header row): qudt:unit -> qudt:unit(separator=",") ---
$ws.Range("C29").Value = 'qudt:unit(separator=",")'

# --- Row 30: AtmosphericParameter -> IgGSubclasses ---
$ws.Range("A30").Value = "incentive-vars:IgGSubclasses"
$ws.Range("B30").Value = "IgG Subclasses"
$ws.Range("C30").Value = "MFI"
$ws.Range("E30").Value = "levels of IgG subclasses to Flu antigens"
$ws.Range("K30").Value = "https://orcid.org/0000-0003-3277-3107"

# --- Row 31: wind_speed -> IgG1-HA1 ---
$ws.Range("A31").Value = "incentive-vars:IgG1-HA1"
$ws.Range("B31").Value = "IgG1-HA1"
$ws.Range("C31").Value = "MFI"
$ws.Range("E31").Value = "levels of immunoglobulin G-1against HA-1 antigen"
$ws.Range("G31").Value = "incentive-vars:IgGSubclasses"
$ws.Range("H31").Value = ""
$ws.Range("K31").Value = "https://orcid.org/0000-0003-3277-3107"

# --- Row 32: (blank) -> IgG2-HA1 ---
$ws.Range("A32").Value = "incentive-vars:IgG2-HA1"
$ws.Range("B32").Value = "IgG2-HA1"
$ws.Range("C32").Value = "MFI"
$ws.Range("E32").Value = "levels of immunoglobulin G-2 against HA-1 antigen"
$ws.Range("G32").Value = "incentive-vars:IgGSubclasses"
$ws.Range("K32").Value = "https://orcid.org/0000-0003-3277-3107"

# --- Row 33: (blank) -> IgG3-HA1 ---
$ws.Range("A33").Value = "incentive-vars:IgG3-HA1"
$ws.Range("B33").Value = "IgG3-HA1"
$ws.Range("C33").Value = "MFI"
$ws.Range("E33").Value = "levels of immunoglobulin G-3 against HA-1 antigen"
$ws.Range("G33").Value = "incentive-vars:IgGSubclasses"
$ws.Range("K33").Value = "https://orcid.org/0000-0003-3277-3107"

# --- Row 34: (blank) -> IgG4-HA1 ---
$ws.Range("A34").Value = "incentive-vars:IgG4-HA1"
$ws.Range("B34").Value = "IgG4-HA1"
$ws.Range("C34").Value = "MFI"
$ws.Range("E34").Value = "levels of immunoglobulin G-4 against HA-1 antigen"
$ws.Range("G34").Value = "incentive-vars:IgGSubclasses"
$ws.Range("K34").Value = "https://orcid.org/0000-0003-3277-3107"

# --- Row 35: (blank) -> SurfaceAntigens-FC ---
$ws.Range("A35").Value = "incentive-vars:SurfaceAntigens-FC"
$ws.Range("B35").Value = "Surface Antigens-FC"
$ws.Range("C35").Value = "MFI,%"
$ws.Range("E35").Value = "level of surface antigenexpression"

# --- Row 36: (blank) -> CD3 ---
$ws.Range("A36").Value = "incentive-vars:CD3"
$ws.Range("B36").Value = "CD3"
$ws.Range("C36").Value = "MFI,%"
$ws.Range("E36").Value = "level of CD3 expression"
$ws.Range("G36").Value = "incentive-vars:SurfaceAntigens-FC"

# --- Row 37: (blank) -> Cytokines-FC ---
$ws.Range("A37").Value = "incentive-vars:Cytokines-FC"
$ws.Range("B37").Value = "Cytokines-FC"
$ws.Range("C37").Value = "MFI,%"
$ws.Range("E37").Value = "level of cytokine expreesion/secretion"

# --- Row 38: (blank) -> IFNg ---
$ws.Range("A38").Value = "incentive-vars:IFNg"
$ws.Range("B38").Value = "IFNg"
$ws.Range("C38").Value = "MFI,%"
$ws.Range("E38").Value = "level of IFNg expression/secretion"
$ws.Range("G38").Value = "incentive-vars:Cytokines"

# --- Row 39: (blank) -> Chemokines-FC ---
$ws.Range("A39").Value = "incentive-vars:Chemokines-FC"
$ws.Range("B39").Value = "Chemokines-FC"
$ws.Range("C39").Value = "MFI,%"
$ws.Range("E39").Value = "level of chemokine expression"

# --- Expand used range to include column U (rows 1-65) and new row 65 ---
$ws.Range("U1:U65").NumberFormat = "General"
$ws.Range("A65:T65").NumberFormat = "General"

# --- New row 65 ---
$ws.Range("A65").Value = "incentive-vars:"

Write-Host ("Final UsedRange: " + $ws.UsedRange.Address())
